# Update "Datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 11:05"

# Row 10 - India: refresh totals
$ws.Cells.Item(10, 2).Value = 191041
$ws.Cells.Item(10, 3).Value = 432
$ws.Cells.Item(10, 4).Value = 91907
$ws.Cells.Item(10, 5).Value = 93721
$ws.Cells.Item(10, 7).Value = 5
$ws.Cells.Item(10, 8).Value = 5413

# Row 24 - Banglades: refresh totals
$ws.Cells.Item(24, 2).Value = 49534
$ws.Cells.Item(24, 3).Value = 2381
$ws.Cells.Item(24, 4).Value = 10597
$ws.Cells.Item(24, 5).Value = 38265
$ws.Cells.Item(24, 7).Value = 22
$ws.Cells.Item(24, 8).Value = 672

# Row 188 - Botsuana: refresh totals
$ws.Cells.Item(188, 2).Value = 38
$ws.Cells.Item(188, 3).Value = 3
$ws.Cells.Item(188, 5).Value = 17

# Rows 192/193 - swap Gambia/Namibia (Namibia now listed first) with refreshed data
$ws.Cells.Item(192, 1).Value = "Namibia"
$ws.Cells.Item(192, 2).Value = 25
$ws.Cells.Item(192, 3).Value = 1
$ws.Cells.Item(192, 4).Value = 16
$ws.Cells.Item(192, 5).Value = 9
$ws.Cells.Item(192, 8).Value = 0

$ws.Cells.Item(193, 1).Value = "Gambia"
$ws.Cells.Item(193, 2).Value = 25
$ws.Cells.Item(193, 4).Value = 20
$ws.Cells.Item(193, 5).Value = 4
$ws.Cells.Item(193, 8).Value = 1

# Rows 201/202 - swap Belice/Santa Lucia (Santa Lucia now listed first) with refreshed data
$ws.Cells.Item(201, 1).Value = "Santa Lucia"
$ws.Cells.Item(201, 4).Value = 18
$ws.Cells.Item(201, 8).Value = 0

$ws.Cells.Item(202, 1).Value = "Belice"
$ws.Cells.Item(202, 4).Value = 16
$ws.Cells.Item(202, 8).Value = 2

# Rows 210/211 - swap Seychelles/Montserrat (Montserrat now listed first) with refreshed data
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Rows 215/216 - swap San Bartolome/Bonaire (Bonaire now listed first)
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216, 1).Value = "San Bartolome"
